# Updates cryptos list (price/volume columns D & E, plus a couple of
# name/link swaps in B & C) to match the latest GitHub Actions scrape.
#
# Note: some new Price values (column D) look like plain decimal numbers
# (e.g. "1.00", "7.57"). Excel's COM layer auto-converts such strings to
# numeric cells, which would lose the original text formatting (the sheet
# stores prices as text, e.g. "63.567.51" with dots as thousand separators).
# To keep those as text we prefix them with a leading apostrophe (Excel's
# "treat as text" marker) and then reset the cell Style back to "Normal" so
# the transient quote-prefix formatting doesn't linger on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '63.567.51'
$ws.Range("E2").Value = '  +4.47%  '

# Row 3
$ws.Range("D3").Value = '3.421.26'
$ws.Range("E3").Value = '  +5.83%  '

# Row 4
$ws.Range("E4").Value = '  +0.01%  '

# Row 5
$ws.Range("D5").Value = '''575.85'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +6.24%  '

# Row 6
$ws.Range("D6").Value = '''156.88'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.42%  '

# Row 7
$ws.Range("D7").Value = '''1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.11%  '

# Row 8
$ws.Range("D8").Value = '3.428.64'
$ws.Range("E8").Value = '  +5.66%  '

# Row 9
$ws.Range("E9").Value = '  +0.72%  '

# Row 10
$ws.Range("D10").Value = '''7.57'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.73%  '

# Row 11
$ws.Range("E11").Value = '  +7.24%  '

# Row 12
$ws.Range("E12").Value = '  -0.24%  '

# Row 13
$ws.Range("D13").Value = '4.016.45'
$ws.Range("E13").Value = '  +6.02%  '

# Row 14
$ws.Range("E14").Value = '  -0.65%  '

# Row 15
$ws.Range("E15").Value = '  +6.83%  '

# Row 16
$ws.Range("D16").Value = '''27.32'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.35%  '

# Row 17
$ws.Range("D17").Value = '63.675.05'
$ws.Range("E17").Value = '  +4.66%  '

# Row 18
$ws.Range("D18").Value = '3.423.51'
$ws.Range("E18").Value = '  +5.96%  '

# Row 19
$ws.Range("E19").Value = '  +1.64%  '

# Row 20
$ws.Range("D20").Value = '''14.28'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +6.78%  '

# Row 21
$ws.Range("D21").Value = '''8.49'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.51%  '

# Row 22
$ws.Range("D22").Value = '''391.39'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.62%  '

# Row 23
$ws.Range("E23").Value = '  -0.29%  '

# Row 24
$ws.Range("D24").Value = '''0.538'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.90%  '

# Row 25
$ws.Range("D25").Value = '''72.02'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.69%  '

# Row 26
$ws.Range("E26").Value = '  +18.78%  '

# Row 27
$ws.Range("D27").Value = '''9.53'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +9.58%  '

# Row 28
$ws.Range("E28").Value = '  +5.08%  '

# Row 29
$ws.Range("E29").Value = '  -0.02%  '

# Row 30
$ws.Range("D30").Value = '''6.71'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +7.92%  '

# Row 31
$ws.Range("D31").Value = '''1.39'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +12.72%  '

# Row 32
$ws.Range("E32").Value = '  +6.62%  '

# Row 33
$ws.Range("D33").Value = '''5.79'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +7.99%  '

# Row 34
$ws.Range("D34").Value = '''23.56'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.20%  '

# Row 36
$ws.Range("E36").Value = '  +3.18%  '

# Row 37
$ws.Range("E37").Value = '  +6.49%  '

# Row 38
$ws.Range("D38").Value = '''158.52'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.11%  '

# Row 39
$ws.Range("D39").Value = '''28.14'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +6.22%  '

# Row 40
$ws.Range("D40").Value = '''0.0783'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +9.48%  '

# Row 41
$ws.Range("E41").Value = '  +8.96%  '

# Row 42
$ws.Range("D42").Value = '2.869.36'
$ws.Range("E42").Value = '  +2.29%  '

# Row 43
$ws.Range("E43").Value = '  +1.93%  '

# Row 44
$ws.Range("B44").Value = 'Mantle'
$ws.Range("C44").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D44").Value = '''0.769'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +5.98%  '

# Row 45
$ws.Range("B45").Value = 'OKB'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D45").Value = '''41.82'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.62%  '

# Row 46
$ws.Range("E46").Value = '  +2.34%  '

# Row 47
$ws.Range("D47").Value = '''1.09'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +9.48%  '

# Row 48
$ws.Range("D48").Value = '3.470.69'
$ws.Range("E48").Value = '  +6.03%  '

# Row 49
$ws.Range("D49").Value = '''22.58'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +7.79%  '

# Row 50
$ws.Range("E50").Value = '  +2.91%  '

# Row 51
$ws.Range("B51").Value = 'dogwifhat'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D51").Value = '''2.09'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +22.24%  '
